# RUSH RACING - POTENTIAL ASSETS.docx edit
#
# Semantic changes applied:
#   1. Remove the "pbr-sand-materials-free-160123" bullet under "Terrain Textures:".
#   2. Remove the "sport-car-3d-model-88076" bullet under "Cars:".
#   3. Re-stamp the document's "last edit" position (the _GoBack bookmark) onto the
#      title, between "RUSH RACING - PO" and "TENTIAL ASSETS:" - this is what Word
#      does automatically whenever the most recent typing/editing happened there.

$d = $word.ActiveDocument

# --- 1 & 2: delete the two asset-link bullet paragraphs -------------------
# Walk paragraphs back-to-front so deleting one never shifts the index of a
# paragraph we still need to examine.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("https://assetstore.unity.com/packages/2d/textures-materials/floors/pbr-sand-materials-free-160123")) {
        $p.Range.Delete()
    }
    elseif ($t.StartsWith("https://assetstore.unity.com/packages/3d/characters/sport-car-3d-model-88076")) {
        $p.Range.Delete()
    }
}

# --- 3: move the _GoBack bookmark into the title --------------------------
# The title is the very first text in the document (before any hyperlink
# field codes), so plain string offsets from Content.Text line up with Word
# range offsets here. Split right after "...RUSH RACING - PO" and before
# "TENTIAL ASSETS:".
$full = $d.Content.Text
$idxPotential = $full.IndexOf("POTENTIAL")
$splitPoint = $idxPotential + 2

$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
